$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.68
$ws.Range("J2").Value = 4
$ws.Range("L2").Value = 1.39
$ws.Range("P2").Value = 1.89
$ws.Range("Q2").Value = 1.95
$ws.Range("W2").Value = 2.46

# --- Row 4 updates ---
$ws.Range("O4").Value = 1.47
$ws.Range("S4").Value = 4.9
$ws.Range("T4").Value = 2.04
$ws.Range("Y4").Value = 7.8

# --- Row 5 updates ---
$ws.Range("F5").Value = 2.8
$ws.Range("G5").Value = 2.82
$ws.Range("H5").Value = 2.92
$ws.Range("I5").Value = 2.94
$ws.Range("L5").Value = 1.5
$ws.Range("P5").Value = 1.72
$ws.Range("U5").Value = 1.98
$ws.Range("W5").Value = 1.54
$ws.Range("AC5").Value = 7
$ws.Range("AJ5").Value = 44

# --- New row 7: Portuguese Segunda Liga match ---
$ws.Range("A7").Value = "Portuguese Segunda Liga"
# "2026-01-15" looks like a real date to Excel's type-inference, so force
# text interpretation, write it, then drop back to the sheet's default
# (unstyled) cell style to mirror the other rows' plain text date cells.
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2026-01-15"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "17:15:00"
$ws.Range("D7").Value = "Vizela"
$ws.Range("E7").Value = "Pacos Ferreira"
$ws.Range("F7").Value = 1.51
$ws.Range("G7").Value = 1.83
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 3.55
$ws.Range("K7").Value = 7.6
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.01
$ws.Range("O7").Value = 1.31
$ws.Range("P7").Value = 1.85
$ws.Range("Q7").Value = 1.8
$ws.Range("R7").Value = 1.1
$ws.Range("S7").Value = 1.01
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.01
$ws.Range("V7").Value = 1.12
$ws.Range("W7").Value = 2.2
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000
